$wb = $excel.ActiveWorkbook

# --- Insert the new "IVY11" worksheet right after "Tickers" (position 2) ---
$indexSheet = $wb.Worksheets.Item("Index")
$ws = $wb.Worksheets.Add($indexSheet)
$ws.Name = "IVY11"

# (look this up only after the sheet collection has settled post-Add)
$tickTest = $wb.Worksheets.Item("TickTest")

# --- Populate header + data rows (A:D) ---
$ws.Range("A1").Value = "Company Name"
$ws.Range("B1").Value = "StockSymbol"
$ws.Range("C1").Value = "Security_Type"
$ws.Range("D1").Value = "Exchange"

$ws.Range("A2").Value = "Real Estate US"
$ws.Range("B2").Value = "VNQ"
$ws.Range("C2").Value = "STK"
$ws.Range("D2").Value = "ISLAND"

$ws.Range("A3").Value = "Real Estate ex-US"
$ws.Range("B3").Value = "VNQI"
$ws.Range("C3").Value = "STK"
$ws.Range("D3").Value = "ISLAND"

$ws.Range("A4").Value = "Bonds US"
$ws.Range("B4").Value = "BND"
$ws.Range("C4").Value = "STK"
$ws.Range("D4").Value = "ISLAND"

$ws.Range("A5").Value = "Bonds ex-US"
$ws.Range("B5").Value = "BNDX"
$ws.Range("C5").Value = "STK"
$ws.Range("D5").Value = "ISLAND"

$ws.Range("A6").Value = "Energy Rohstoffe XLE besser"
$ws.Range("B6").Value = "DBE"
$ws.Range("C6").Value = "STK"
$ws.Range("D6").Value = "ISLAND"

$ws.Range("A7").Value = "Edelmetalle"
$ws.Range("B7").Value = "DBP"
$ws.Range("C7").Value = "STK"
$ws.Range("D7").Value = "ISLAND"

$ws.Range("A8").Value = "Industriemetalle"
$ws.Range("B8").Value = "DBB"
$ws.Range("C8").Value = "STK"
$ws.Range("D8").Value = "ISLAND"

$ws.Range("A9").Value = "Agriculture"
$ws.Range("B9").Value = "DBA"
$ws.Range("C9").Value = "STK"
$ws.Range("D9").Value = "ISLAND"

$ws.Range("A10").Value = "Emergingmarket"
$ws.Range("B10").Value = "EEM"
$ws.Range("C10").Value = "STK"
$ws.Range("D10").Value = "ISLAND"

$ws.Range("A11").Value = "Stocks nonUS"
$ws.Range("B11").Value = "EFA"
$ws.Range("C11").Value = "STK"
$ws.Range("D11").Value = "ISLAND"

$ws.Range("A12").Value = "Stocks US"
$ws.Range("B12").Value = "VTI"
$ws.Range("C12").Value = "STK"
$ws.Range("D12").Value = "ISLAND"

$ws.Range("A13").Value = "Energieaktien Öl"
$ws.Range("B13").Value = "XLE"
$ws.Range("C13").Value = "STK"
$ws.Range("D13").Value = "ISLAND"

$ws.Range("A14").Value = "Cash 1-3 y US-treasury"
$ws.Range("B14").Value = "SHY"
$ws.Range("C14").Value = "STK"
$ws.Range("D14").Value = "ISLAND"

# --- Match header-cell format in C1 (same style used on the other sheets) ---
$tickTest.Range("C1").Copy()
$ws.Range("C1").PasteSpecial(-4122)

# --- Column A width ---
$ws.Columns.Item(1).ColumnWidth = 19.619791666666668

# --- This new sheet becomes the active tab / selected sheet ---
$ws.Activate()
$ws.Range("B17").Select()
